$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 76.03
$ws.Range("C2").Value = 'The applicant has a strong background in front-end development, with experience in ReactJS and JavaScript, which aligns well with the requirements of the company. The candidate has worked on projects involving ReactJS, Django, Python, JavaScript, and various data analysis and computer vision skills, showcasing a diverse skill set. However, the applicant lacks proficiency in CSS and NextJS, which are also required by the company. The projects demonstrate the candidate''s capability to work on user-friendly UI development and data analysis pipelines, providing valuable experience relevant to the job role. The applicant''s technical skills and project experience make them a suitable candidate for the position, although further proficiency in CSS and NextJS would be beneficial.'
$ws.Range("E2").Value = 'The applicant exhibits strong communication skills, evident from the responses provided during the interview. The applicant also showcases a high level of adaptability and a willingness to learn from Japanese work culture, indicating a positive attitude towards the prospect of working in Japan. Furthermore, the applicant''s expressed interest in immersing in Japanese culture and values, along with the desire for mentorship and skill development, demonstrates a strong commitment to personal and professional growth. Overall, the applicant''s responses reflect a well-rounded personality and a genuine willingness to embrace the challenges and opportunities of working in Japan.'

# Row 3 updates
$ws.Range("B3").Value = 90.81
$ws.Range("C3").Value = 'The candidate possesses strong technical skills in frontend development using ReactJS and has also showcased proficiency in technologies such as Flutter, WebRTC, NodeJS, ExpressJS, and HTML which could be beneficial in creating interactive and real-time user interfaces. The candidate has demonstrated a strong ability to work in a team and has shown adaptability by working with different technology stacks in various projects. However, the candidate lacks direct experience with NextJS, a skill specifically mentioned in the job description, which might indicate a gap in full compliance with the required skills. Nevertheless, the candidate''s diverse skill set, particularly in user interface development and real-time communication, aligns well with the company''s needs.'
$ws.Range("E3").Value = 'The applicant has shown a strong willingness to work in Japan through their expressed interest in the country''s AI advancements, work culture, and learning a new language. They have also demonstrated their adaptability to new cultures and a desire to work in a team environment. However, their preference to work in India or the US as well might indicate some flexibility, albeit the overall positive attitude towards Japan is evident. Additionally, their commitment to continuous learning and improvement, along with their ability to function well within a team, supports a high score for their suitability to work in Japan.'

# Row 4 updates
$ws.Range("B4").Value = 77.73
$ws.Range("C4").Value = 'The applicant has a strong background in backend development and machine learning, with experience in using Django Rest Framework, PyTorch, Tensorflow, and Keras. While the job description primarily focuses on frontend development, the applicant''s experience in ReactJS and Django Rest Framework demonstrates adaptability and the ability to grasp new concepts. The applicant has also showcased skills in image processing and computer vision, which may not directly align with the job requirements but demonstrate a breadth of technical knowledge. Additionally, the applicant has experience in developing a library management system using Django Rest Framework and ReactJS, which showcases proficiency in web development. Despite lacking experience in specific technologies like NextJS, the applicant''s diverse technical skills and willingness to learn make them a strong candidate for the position.'
$ws.Range("E4").Value = 'The applicant has displayed excellent communication skills, adaptability, and a strong enthusiasm for learning and contributing to cutting-edge technologies. They have shown a willingness to work in Japan and are open to adapting to a new culture. The applicant''s collaborative problem-solving approach and ability to work well in a team align with the company''s values, making them a suitable candidate for the job and a potential fit for the company''s work environment in Japan.'

# New row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 72.7
$ws.Range("C5").Value = 'The applicant has a good understanding of front-end development, as demonstrated by the projects ''Video Conferencing Project'', ''TEDxIITGuwahati Website'', and ''Alcheringa Pass Portal''. These projects required the use of JavaScript, HTML, and CSS, which are essential for the role of a Frontend Engineer Intern. The applicant''s experience in API integration also showcases their ability to work with different technologies and integrate them into web applications. However, the applicant lacks direct experience with ReactJS and NextJS, which are specifically mentioned in the job description. While they possess the foundational skills required, the absence of experience with these specific technologies may impact their immediate effectiveness in the role.'
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 'The applicant''s answers demonstrate strong communication skills, adaptability, and a willingness to embrace new challenges and cultures. They have expressed a genuine interest in Japanese work culture and are open to working in Japan after graduation, showcasing their willingness to relocate. Additionally, their acknowledgment of areas for improvement and their proactive steps to address them reflect a growth-oriented mindset. Overall, the applicant exhibits qualities that align with the company''s soft skills requirements and a positive attitude towards working in Japan.'

